{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: insert a new \"Meta description\" paragraph right after the\n//           Heading1 title paragraph at the top of the document. The\n//           paragraph has two runs: \"Meta description\" (bold) followed\n//           by \": Read our review of Dragon Egg, a free slot game\n//           featuring dragons and symbols of wealth, with a\n//           straightforward interface and medium volatility.\" (not bold).\n//\n// Change 2: near the end of the document, remove the bold\n//           \"Play Dragon Egg Free Slot Game - Review 2021\" paragraph\n//           entirely, and replace the text of the following italic\n//           paragraph with the DALLE image prompt (keeping the italic\n//           formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---- Change 1: add the Meta description paragraph under the H1 title ----\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.styleBuiltIn = Word.Style.normal;\nawait context.sync();\n\nconst afterLabelText =\n  \": Read our review of Dragon Egg, a free slot game featuring dragons and symbols of wealth, with a straightforward interface and medium volatility.\";\nmetaPara.insertText(afterLabelText, \"End\");\nawait context.sync();\n\nconst labelRange = metaPara.insertText(\"Meta description\", \"Start\");\nlabelRange.font.bold = true;\nawait context.sync();\n\n// ---- Change 2: drop the bold title paragraph, repurpose the italic one ----\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nconst boldTitlePara = paragraphs.items[count - 2];\nconst italicDescPara = paragraphs.items[count - 1];\n\nconst dallePrompt =\n  \"Prompt for DALLE: Create a feature image for Dragon Egg, a slot game by Tom Horn, in a cartoon style. The image should feature a happy Maya warrior wearing glasses. The warrior should have a confident expression on their face and be holding a golden dragon egg in one hand, as if they have just won it in the slot game. The background should be a dark cave, with shadows of dragons visible in the background. The image should be eye-catching and convey the excitement of winning big in the game.\";\n\nboldTitlePara.delete();\nconst descRange = italicDescPara.getRange();\ndescRange.insertText(dallePrompt, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Change 1: insert a new \"Meta description\" paragraph right after the\n#           Heading1 title paragraph at the top of the document. The\n#           paragraph has two runs: \"Meta description\" (bold) followed\n#           by \": Read our review of Dragon Egg, a free slot game\n#           featuring dragons and symbols of wealth, with a\n#           straightforward interface and medium volatility.\" (not bold).\n#\n# Change 2: near the end of the document, remove the bold\n#           \"Play Dragon Egg Free Slot Game - Review 2021\" paragraph\n#           entirely, and replace the text of the following italic\n#           paragraph with the DALLE image prompt (keeping the italic\n#           formatting).\n\n$d = $word.ActiveDocument\n\n# ---- Change 1: add the Meta description paragraph under the H1 title ----\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Style = \"Normal\"\n$metaRange = $metaPara.Range\n$metaRange.Text = \"Meta description: Read our review of Dragon Egg, a free slot game featuring dragons and symbols of wealth, with a straightforward interface and medium volatility.\"\n\n$label = \"Meta description\"\n$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $label.Length)\n$boldRange.Font.Bold = $true\n\n# ---- Change 2: drop the bold title paragraph, repurpose the italic one ----\n$count = $d.Paragraphs.Count\n$boldTitlePara = $d.Paragraphs($count - 1)\n$btr = $boldTitlePara.Range\n$d.Range($btr.Start, $btr.End).Delete()\n\n$count2 = $d.Paragraphs.Count\n$italicDescPara = $d.Paragraphs($count2)\n$idr = $italicDescPara.Range\n$dallePrompt = \"Prompt for DALLE: Create a feature image for Dragon Egg, a slot game by Tom Horn, in a cartoon style. The image should feature a happy Maya warrior wearing glasses. The warrior should have a confident expression on their face and be holding a golden dragon egg in one hand, as if they have just won it in the slot game. The background should be a dark cave, with shadows of dragons visible in the background. The image should be eye-catching and convey the excitement of winning big in the game.\"\n$d.Range($idr.Start, $idr.End).Text = $dallePrompt\n"}
